$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two header cells (C1/D1): "P. BOLSA" -> "BOLSA", "P. MUESTRA" -> "MUESTRA"
# (set D1 first so "MUESTRA" is appended to the shared-strings table before "BOLSA",
# matching the order produced by the original edit)
$ws.Range("D1").Value = "MUESTRA"
$ws.Range("C1").Value = "BOLSA"

# Row 1 no longer needs the taller wrapped height now that the headers are shorter
$ws.Rows("1").RowHeight = 15.75

# Update the current selection to match the diff
$ws.Range("C1:C2").Select()
